$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("D3").Value = 44585
$ws.Range("M3").Value = 160
$ws.Range("N3").Value = 6500
$ws.Range("O3").Value = 7000
$ws.Range("P3").Value = 6750
$ws.Range("R3").Value = "Provincia de Curicó"
$ws.Range("S3").Value = 3375

# Row 4
$ws.Range("D4").Value = 44606
$ws.Range("M4").Value = 45
$ws.Range("N4").Value = 7000
$ws.Range("O4").Value = 7000
$ws.Range("P4").Value = 7000
$ws.Range("R4").Value = "Provincia de Linares"
$ws.Range("S4").Value = 3500

# Row 5
$ws.Range("D5").Value = 44614
$ws.Range("M5").Value = 45
$ws.Range("N5").Value = 6000
$ws.Range("O5").Value = 6000
$ws.Range("P5").Value = 6000
$ws.Range("S5").Value = 3000

# Row 6
$ws.Range("D6").Value = 44627
$ws.Range("M6").Value = 45
$ws.Range("R6").Value = "Provincia de Linares"

# Row 7
$ws.Range("D7").Value = 44589
$ws.Range("M7").Value = 60
$ws.Range("N7").Value = 6000
$ws.Range("O7").Value = 6000
$ws.Range("P7").Value = 6000
$ws.Range("R7").Value = "Provincia de Curicó"
$ws.Range("S7").Value = 3000

# Row 8
$ws.Range("D8").Value = 44209
$ws.Range("M8").Value = 58
$ws.Range("N8").Value = 6000
$ws.Range("O8").Value = 6000
$ws.Range("P8").Value = 6000
$ws.Range("S8").Value = 3000

# Row 9
$ws.Range("D9").Value = 44587
$ws.Range("M9").Value = 165
$ws.Range("N9").Value = 6500
$ws.Range("O9").Value = 7000
$ws.Range("P9").Value = 6742
$ws.Range("S9").Value = 3371

# Row 10
$ws.Range("D10").Value = 44214
$ws.Range("M10").Value = 48

# Row 12
$ws.Range("D12").Value = 44586
$ws.Range("M12").Value = 80
$ws.Range("N12").Value = 7000
$ws.Range("P12").Value = 7000
$ws.Range("S12").Value = 3500

# Row 14
$ws.Range("D14").Value = 44592
$ws.Range("M14").Value = 30
$ws.Range("N14").Value = 8000
$ws.Range("O14").Value = 8000
$ws.Range("P14").Value = 8000
$ws.Range("R14").Value = "Provincia de Linares"
$ws.Range("S14").Value = 4000

# Row 15
$ws.Range("D15").Value = 44211
$ws.Range("M15").Value = 45
$ws.Range("N15").Value = 6000
$ws.Range("O15").Value = 6000
$ws.Range("P15").Value = 6000
$ws.Range("R15").Value = "Provincia de Curicó"
$ws.Range("S15").Value = 3000
